$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Fecha(D), Variedad(H), Volumen(J), PrecioMin(K), PrecioMax(L), PrecioProm(M), Unidad(N), Origen(O), PrecioKg(P)
$rows = @(
    @(2, 44595, "Perfection", 50, 26000, 28000, 27200, "$/saco 25 kilos", "Carahue", 1088),
    @(3, 44615, "Sin especificar", 200, 28000, 30000, 29000, "$/saco 25 kilos", "Carahue", 1160),
    @(4, 44518, "Perfection", 350, 14000, 15000, 14571, "$/saco 25 kilos", "Región del Maule", 583),
    @(5, 44532, "Sin especificar", 250, 14000, 15000, 14400, "$/saco 25 kilos", "Región del Maule", 576),
    @(6, 44503, "Perfection", 200, 15000, 16000, 15500, "$/malla 25 kilos", "Provincia de Limarí", 620),
    @(7, 44539, "Sin especificar", 50, 13000, 14000, 13400, "$/saco 25 kilos", "Región del Maule", 536),
    @(8, 44342, "Perfection", 60, 30000, 32000, 31000, "$/malla 25 kilos", "Provincia de Limarí", 1240),
    @(9, 44505, "Perfection", 210, 6500, 7000, 6714, "$/malla 25 kilos", "Región del Maule", 269),
    @(10, 44483, "Perfection", 220, 19000, 20000, 19455, "$/saco 25 kilos", "Región Metropolitana", 778),
    @(11, 44643, "Perfection", 90, 25000, 26000, 25444, "$/saco 25 kilos", "Carahue", 1018),
    @(12, 44673, "Sin especificar", 220, 25000, 26000, 25455, "$/saco 25 kilos", "Carahue", 1018),
    @(13, 44540, "Sin especificar", 110, 16000, 17000, 16545, "$/saco 25 kilos", "Región del Maule", 662),
    @(14, 44659, "Sin especificar", 140, 24000, 25000, 24571, "$/saco 25 kilos", "Carahue", 983),
    @(15, 44671, "Perfection", 110, 25000, 26000, 25545, "$/saco 25 kilos", "Carahue", 1022),
    @(16, 44533, "Perfection", 80, 14000, 15000, 14375, "$/malla 25 kilos", "Región del Maule", 575),
    @(17, 44162, "Sin especificar", 100, 17000, 18000, 17500, "$/saco 25 kilos", "Región del Maule", 700),
    @(18, 44399, "Perfection", 50, 39000, 40000, 39600, "$/malla 25 kilos", "Provincia de Huasco", 1584),
    @(19, 44545, "Perfection", 180, 15000, 16000, 15444, "$/saco 25 kilos", "Carahue", 618),
    @(20, 44631, "Perfection", 150, 24000, 25000, 24467, "$/saco 25 kilos", "Carahue", 979),
    @(21, 44519, "Perfection", 240, 17000, 18000, 17583, "$/saco 25 kilos", "Carahue", 703),
    @(22, 44657, "Sin especificar", 250, 24000, 25000, 24400, "$/saco 25 kilos", "Carahue", 976),
    @(23, 44335, "Perfection", 100, 30000, 32000, 31000, "$/malla 25 kilos", "Provincia de Huasco", 1240),
    @(24, 44589, "Perfection", 160, 22000, 23000, 22500, "$/malla 25 kilos", "Carahue", 900),
    @(25, 44629, "Perfection", 35, 25000, 26000, 25429, "$/saco 25 kilos", "Región Metropolitana", 1017),
    @(26, 44454, "Perfection", 100, 36000, 38000, 37000, "$/malla 25 kilos", "Provincia de Limarí", 1480),
    @(27, 44517, "Perfection", 110, 17000, 18000, 17455, "$/saco 25 kilos", "Región del Maule", 698),
    @(28, 44482, "Perfection", 130, 24000, 25000, 24385, "$/saco 25 kilos", "Región de O'Higgins", 975),
    @(29, 44328, "Perfection", 100, 33000, 34000, 33500, "$/malla 25 kilos", "Provincia de Huasco", 1340),
    @(30, 44496, "Perfection", 250, 14000, 15000, 14520, "$/malla 25 kilos", "Provincia de Huasco", 581),
    @(31, 44512, "Perfection", 100, 14000, 15000, 14500, "$/saco 25 kilos", "Región del Maule", 580)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D - Fecha
    $ws.Cells.Item($r, 8).Value = $row[2]   # H - Variedad
    $ws.Cells.Item($r, 10).Value = $row[3]  # J - Volumen
    $ws.Cells.Item($r, 11).Value = $row[4]  # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[5]  # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[6]  # M - Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $row[7]  # N - Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $row[8]  # O - Origen
    $ws.Cells.Item($r, 16).Value = $row[9]  # P - Precio $/Kg
}
